$d = $word.ActiveDocument

# 1. Remove the stale _GoBack bookmark left after "View User allocation"
#    (Word leaves a fresh one at the location of the last edit, below.)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. "View upcoming projects" -> "View projects"
$heading = $d.Content
[void]$heading.Find.Execute("View upcoming projects")
$headingStart = $heading.Start

# Delete just the "upcoming " word (keeps "View " / "projects" runs intact
# conceptually) using offsets relative to the unique match found above.
$toDelete = $d.Range($headingStart + 5, $headingStart + 14)
$toDelete.Text = ""

# The deletion merges the now-adjacent, identically formatted runs into a
# single run, swallowing the separate leading-space run that precedes
# "View projects". Re-split them back apart (toggling a character
# property forces Word to re-break the run without altering the text or
# leaving a residual formatting override behind).
$reSplit = $d.Content
[void]$reSplit.Find.Execute("View projects")
$reSplit.Bold = $true
$reSplit.Bold = $false

# 3. "Add Resource Requirement to project" -> same text, but the edit
#    point lands after "Add Resource Req", which is where Word leaves
#    the _GoBack bookmark marking the last editing location.
$target = $d.Content
[void]$target.Find.Execute("Add Resource Requirement to project")
$splitPoint = $target.Start + 16
$bmRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bmRange)
